$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Random Forest Results Test Data")
$ws1.Range("A2").Value = "Degree_Rede_EntreajudaLabur"
$ws1.Range("B2").Value = 0.5388741376146789
$ws1.Range("C2").Value = 3.856245184512639
$ws1.Range("D2").Value = 2.851153846153846
$ws1.Range("A3").Value = "OutDeg_Var.Dep_RedeControlExtAusencia"
$ws1.Range("B3").Value = 0.4260617629889084
$ws1.Range("C3").Value = 1.205973912206677
$ws1.Range("D3").Value = 0.9973076923076923
$ws1.Range("A4").Value = "OutDeg_Var.Dep_ApoioSpecDiqCint"
$ws1.Range("B4").Value = 0.6320215005318817
$ws1.Range("C4").Value = 2.949441928960069
$ws1.Range("D4").Value = 1.446923076923077
$ws1.Range("A5").Value = "OutDeg_Var.Dep_ApoioSpecTubos"
$ws1.Range("B5").Value = 0.1872776684330054
$ws1.Range("C5").Value = 15.74024386479056
$ws1.Range("D5").Value = 3.778461538461538
$ws1.Range("A6").Value = "OutDeg_Var.Dep_RedApoiLevntDiCin"
$ws1.Range("B6").Value = 0.094672463768116
$ws1.Range("C6").Value = 0.9119590071758877
$ws1.Range("D6").Value = 0.6038461538461538
$ws1.Range("A7").Value = "OutDeg_RedCont_DiqPriqDiqCint"
$ws1.Range("B7").Value = 0.7113870758053724
$ws1.Range("C7").Value = 2.113348600608122
$ws1.Range("D7").Value = 1.515
$ws1.Range("A8").Value = "OutDeg_RedeRepar_DiqPriqDiqCint"
$ws1.Range("B8").Value = 0.4282970607734806
$ws1.Range("C8").Value = 1.95623559223006
$ws1.Range("D8").Value = 1.428076923076923
$ws1.Range("A9").Value = "OutDeg_Red-Val_Cont1Rep2DiqCin"
$ws1.Range("B9").Value = 0.0514886544671691
$ws1.Range("C9").Value = 2.552940988802462
$ws1.Range("D9").Value = 1.673076923076923
$ws1.Range("A10").Value = "OutDeg_Red-Val_Cont1Rep2DiqPrq"
$ws1.Range("B10").Value = 0.584054759383457
$ws1.Range("C10").Value = 3.960858007281522
$ws1.Range("D10").Value = 2.860384615384616
$ws1.Range("A11").Value = "OutDeg_ContRep2_DiqPrqCin"
$ws1.Range("B11").Value = 0.6012632724446172
$ws1.Range("C11").Value = 4.267566411542189
$ws1.Range("D11").Value = 3.155384615384615
$ws1.Range("A12").Value = "OutDeg_Var.Dep_RedePartilhaAgu"
$ws1.Range("B12").Value = -0.1793870833333333
$ws1.Range("C12").Value = 0.9151145367573479
$ws1.Range("D12").Value = 0.655

$ws2 = $wb.Worksheets.Item("Random Forest Results All Data")
$ws2.Range("A2").Value = "Degree_Rede_EntreajudaLabur"
$ws2.Range("B2").Value = 0.7983547168849485
$ws2.Range("C2").Value = 2.388480104439291
$ws2.Range("D2").Value = 1.612558139534884
$ws2.Range("A3").Value = "OutDeg_Var.Dep_RedeControlExtAusencia"
$ws2.Range("B3").Value = 0.7961400953836425
$ws2.Range("C3").Value = 0.6818573585174109
$ws2.Range("D3").Value = 0.5025581395348837
$ws2.Range("A4").Value = "OutDeg_Var.Dep_ApoioSpecDiqCint"
$ws2.Range("B4").Value = 0.8231424815978156
$ws2.Range("C4").Value = 1.492105582647368
$ws2.Range("D4").Value = 0.6656589147286822
$ws2.Range("A5").Value = "OutDeg_Var.Dep_ApoioSpecTubos"
$ws2.Range("B5").Value = 0.3459148404748377
$ws2.Range("C5").Value = 7.251202733703652
$ws2.Range("D5").Value = 1.21968992248062
$ws2.Range("A6").Value = "OutDeg_Var.Dep_RedApoiLevntDiCin"
$ws2.Range("B6").Value = 0.7933523227712138
$ws2.Range("C6").Value = 0.5267517839176865
$ws2.Range("D6").Value = 0.3326356589147287
$ws2.Range("A7").Value = "OutDeg_RedCont_DiqPriqDiqCint"
$ws2.Range("B7").Value = 0.888971319317113
$ws2.Range("C7").Value = 1.18514775188516
$ws2.Range("D7").Value = 0.6989147286821706
$ws2.Range("A8").Value = "OutDeg_RedeRepar_DiqPriqDiqCint"
$ws2.Range("B8").Value = 0.8652813553062826
$ws2.Range("C8").Value = 1.025791432024807
$ws2.Range("D8").Value = 0.6034108527131782
$ws2.Range("A9").Value = "OutDeg_Red-Val_Cont1Rep2DiqCin"
$ws2.Range("B9").Value = 0.8073187668001536
$ws2.Range("C9").Value = 1.29944054747961
$ws2.Range("D9").Value = 0.6992248062015503
$ws2.Range("A10").Value = "OutDeg_Red-Val_Cont1Rep2DiqPrq"
$ws2.Range("B10").Value = 0.864793267161582
$ws2.Range("C10").Value = 2.175913601225931
$ws2.Range("D10").Value = 1.285271317829457
$ws2.Range("A11").Value = "OutDeg_ContRep2_DiqPrqCin"
$ws2.Range("B11").Value = 0.8791652726542241
$ws2.Range("C11").Value = 2.273534965710903
$ws2.Range("D11").Value = 1.352403100775194
$ws2.Range("A12").Value = "OutDeg_Var.Dep_RedePartilhaAgu"
$ws2.Range("B12").Value = 0.5994635493372606
$ws2.Range("C12").Value = 0.4783328606539544
$ws2.Range("D12").Value = 0.2893798449612403
